$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells in this sheet are plain, unstyled inline strings (no explicit
# cell style / number format). New values are prefixed with a leading apostrophe
# so Excel treats them as literal text (preventing numeric-looking strings such
# as "0.9978" or "21.760.01" from being reinterpreted as numbers on assignment).
# The apostrophe prefix makes Excel auto-apply a "quote prefix" text style to the
# cell, so we immediately reset the cell style back to "Normal" to match the
# original (unstyled) cells exactly.

# Row 2
$ws.Range("D2").Value = "'21.713.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.78%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.574.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +6.78%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'0.9978"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -1.04%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.9750"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.87%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'284.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.75%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.3682"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.87%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.3264"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +6.77%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  +7.35%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'41.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.68%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07048"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +6.55%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.9944"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.79%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'20.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +10.82%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'5.792"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +6.00%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'6.478"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.89%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.00001069"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.98%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = "'Dai"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.9728"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.84%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'1.561.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.92%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.06183"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.82%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'73.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.66%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'15.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +10.00%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'5.830"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +6.66%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'11.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +5.00%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'21.707.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +5.44%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'2.325"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.14%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'2.402"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +12.80%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'148.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +5.22%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'18.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.64%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'1.736.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +6.38%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'120.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.71%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'4.053"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.82%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'0.9041"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +10.26%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'5.380"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +8.18%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'0.08160"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.61%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'1.570"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.46%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'5.072"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +7.35%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("B37").Value = "'Aptos"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'11.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +10.21%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("B38").Value = "'TrustWalletToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1.235"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.24%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.06011"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.38%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'8.123"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +7.05%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  +6.14%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  +6.78%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.9710"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.44%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.5701"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +8.01%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'12.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +7.01%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'3.602"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.84%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.5619"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +8.30%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'124.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.83%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'1.917"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +7.92%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.06722"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.06%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'71.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +6.10%  "
$ws.Range("E51").Style = "Normal"
